# "Se trabaja en el formulario de nuevas tareas"
#
# - Column K ("Comentarios") status notes updated:
#     * old note "Pendiente el flujo de ocultarse el formulario"      -> "Listo"
#     * old note "Pendiente  botones de olvido y regristrar el usuario" -> "Listo, falta olvido contraseña"
#   (these are shared strings already referenced by K2/K3)
# - Several additional requirement rows get marked "Listo" in column K
#   (rows 6, 7, 20, 21, 22, 23, 24, 25, 26)
# - Header row highlight fill removed from A2/A3 (now matches the plain
#   vertical-centered style used elsewhere, e.g. A4)
# - The view's scroll/selection moved down a bit (selection -> K8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two existing status comments (reused by K2 / K3) ---
$ws.Range("K2").Value = "Listo"
$ws.Range("K3").Value = "Listo, falta olvido contraseña"

# --- Mark additional requirement rows as "Listo" in column K ---
$ws.Range("K6").Value = "Listo"
$ws.Range("K7").Value = "Listo"
$ws.Range("K20").Value = "Listo"
$ws.Range("K21").Value = "Listo"
$ws.Range("K22").Value = "Listo"
$ws.Range("K23").Value = "Listo"
$ws.Range("K24").Value = "Listo"
$ws.Range("K25").Value = "Listo"
$ws.Range("K26").Value = "Listo"

# --- Remove the highlight fill from A2/A3 so they match the plain style ---
# (copy the formatting of A4, which already uses the no-fill / vertical-center style)
$ws.Range("A4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Move the active selection in the sheet view ---
$ws.Activate()
$ws.Range("K8").Select()
